$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1) Metadata sheet: bump the "Date" value (row 8, column B) to the new
#    publication timestamp.
# ---------------------------------------------------------------------------
$meta = $wb.Worksheets.Item("Metadata")
$meta.Range("B8").Value = "2024-03-22T16:25:12+00:00"

# ---------------------------------------------------------------------------
# 2) Elements sheet: the two "Mapping" columns (AK = col 37, AL = col 38)
#    swap places - the French "business mapping" column now comes first
#    (AK), and "RIM Mapping" moves to the second slot (AL). That means both
#    the header text and every data row's contents need to trade places
#    between the two columns, and the column widths need to follow the
#    content (the wider French label now sits in AK, the narrower "RIM
#    Mapping" in AL).
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Elements")

$lastRow = $ws.UsedRange.Rows.Count
if ($lastRow -lt 6) { $lastRow = 6 }

for ($r = 1; $r -le $lastRow; $r++) {
    $akCell = $ws.Range("AK$r")
    $alCell = $ws.Range("AL$r")
    $akVal = $akCell.Value()
    $alVal = $alCell.Value()
    # Only touch rows where the two columns actually differ - both columns
    # are blank on several rows, and re-writing "" back onto an already
    # empty cell needlessly perturbs its stored representation.
    if ("$akVal" -ne "$alVal") {
        $akCell.Value = $alVal
        $alCell.Value = $akVal
    }
}

# Swap the column widths to match (AK becomes the wide "Spécification
# métier" column, AL becomes the narrower "RIM Mapping" column).
$ws.Columns.Item(37).ColumnWidth = 66.83
$ws.Columns.Item(38).ColumnWidth = 24.165
